# "ajout de log dans le test connexion"
# Update the connection-test log sheet (List_ID):
#   - F3 (exp)    : fail -> pass
#   - G3 (Status) : PASS -> FAIL
#   - append 5 new log rows (11-15) with date_naissance / email / mdp / exp
#   - move the active selection to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_ID")

# --- fix up the existing row 3 result ---------------------------------
$ws.Range("F3").Value = "pass"
$ws.Range("G3").Value = "FAIL"

# Helper: write a literal text value into a cell. Some of the
# "date_naissance" strings below (DD/MM/YYYY) look like valid dates to
# Excel's auto-detection and would otherwise be silently turned into a
# date serial number - so re-apply as explicit text whenever that happens.
function Set-TextValue($cell, $text) {
    $cell.Value = $text
    $current = $cell.Value()
    if ($current -is [DateTime]) {
        $cell.NumberFormat = "@"
        $cell.Value = $text
    }
}

# --- new log rows -------------------------------------------------------
$newRows = @(
    @("12/10/1978", "ujpshrhgrcukjbb@gmail.com", "kseecZQOLJ5", "pass"),
    @("30/12/1989", "ajcsxgxcjbnrypl@gmail.com", "diwmcWKCOB5", "pass"),
    @("04/09/1977", "euygnnbmxacyhco@gmail.com", "sstglAVTKQ5", "pass"),
    @("14/12/1981", "erobvegrhbxuncq@gmail.com", "alqmtSAKWZ5", "pass"),
    @("23/11/2004", "orkmjqjqnosdhqi@gmail.com", "sgshzWELWB5", "pass")
)

$row = 11
foreach ($entry in $newRows) {
    Set-TextValue $ws.Cells.Item($row, 3) $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]
    $ws.Cells.Item($row, 5).Value = $entry[2]
    $ws.Cells.Item($row, 6).Value = $entry[3]

    $row = $row + 1
}

# --- move the selection, matching the recorded cursor position ---------
$ws.Range("G11").Select()
